$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028368881697471
$ws.Range("D2").Value = 1.03675252746226
$ws.Range("E2").Value = 1.028335307673828
$ws.Range("F2").Value = 1.043700748209906
$ws.Range("I2").Value = 1.032044799016489
$ws.Range("J2").Value = 1.033521584377905
$ws.Range("K2").Value = 1.039545370101956
$ws.Range("L2").Value = 1.031152423627581
$ws.Range("M2").Value = 1.046473873974521
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029373999598164
$ws.Range("D3").Value = 1.03753676966259
$ws.Range("E3").Value = 1.029190568526886
$ws.Range("F3").Value = 1.044677412616198
$ws.Range("I3").Value = 1.032202851175256
$ws.Range("J3").Value = 1.034166996776917
$ws.Range("K3").Value = 1.040139346069516
$ws.Range("L3").Value = 1.031815470601121
$ws.Range("M3").Value = 1.047261194033815
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030024806538998
$ws.Range("D4").Value = 1.038044376230691
$ws.Range("E4").Value = 1.029744717826858
$ws.Range("F4").Value = 1.045310023185287
$ws.Range("I4").Value = 1.032303792312888
$ws.Range("J4").Value = 1.034584474642875
$ws.Range("K4").Value = 1.040523181880791
$ws.Range("L4").Value = 1.032244605102927
$ws.Range("M4").Value = 1.047770669563352
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03029850779921
$ws.Range("D5").Value = 1.038257808648254
$ws.Range("E5").Value = 1.0299778578884
$ws.Range("F5").Value = 1.04557612562468
$ws.Range("I5").Value = 1.032345909297255
$ws.Range("J5").Value = 1.034759946366277
$ws.Range("K5").Value = 1.04068442435848
$ws.Range("L5").Value = 1.03242503607676
$ws.Range("M5").Value = 1.047984858531301
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030344469411926
$ws.Range("D6").Value = 1.03829364688662
$ws.Range("E6").Value = 1.030017013402576
$ws.Range("F6").Value = 1.04562081433344
$ws.Range("I6").Value = 1.032352962223883
$ws.Range("J6").Value = 1.03478940670953
$ws.Range("K6").Value = 1.040711490494733
$ws.Range("L6").Value = 1.032455332546184
$ws.Range("M6").Value = 1.048020822093474
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030028463348434
$ws.Range("D7").Value = 1.038047227990538
$ws.Range("E7").Value = 1.029747832366743
$ws.Range("F7").Value = 1.045313578259438
$ws.Range("I7").Value = 1.032304356334826
$ws.Range("J7").Value = 1.034586819446918
$ws.Range("K7").Value = 1.040525336891529
$ws.Range("L7").Value = 1.032247015941753
$ws.Range("M7").Value = 1.047773531546597
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028708476953156
$ws.Range("D8").Value = 1.037017534151437
$ws.Range("E8").Value = 1.028624193565721
$ws.Range("F8").Value = 1.044030682654666
$ws.Range("I8").Value = 1.032098488420482
$ws.Range("J8").Value = 1.033739734528176
$ws.Range("K8").Value = 1.039746211379948
$ws.Range("L8").Value = 1.03137648251335
$ws.Range("M8").Value = 1.046739946345816
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026385797786912
$ws.Range("D9").Value = 1.035204277350283
$ws.Range("E9").Value = 1.026649900085692
$ws.Range("F9").Value = 1.041775031529055
$ws.Range("I9").Value = 1.03172556330165
$ws.Range("J9").Value = 1.032245969044679
$ws.Range("K9").Value = 1.038369453064372
$ws.Range("L9").Value = 1.029843283843143
$ws.Range("M9").Value = 1.044918884786762
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024839597962745
$ws.Range("D10").Value = 1.033996310673444
$ws.Range("E10").Value = 1.025337599971992
$ws.Range("F10").Value = 1.040274671176939
$ws.Range("I10").Value = 1.031470142858967
$ws.Range("J10").Value = 1.031249433578329
$ws.Range("K10").Value = 1.037449083502132
$ws.Range("L10").Value = 1.028821731659435
$ws.Range("M10").Value = 1.04370506605795
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024170616518218
$ws.Range("D11").Value = 1.033473468291537
$ws.Range("E11").Value = 1.024770296240013
$ws.Range("F11").Value = 1.039625819357558
$ws.Range("I11").Value = 1.031357934775588
$ws.Range("J11").Value = 1.030817767428594
$ws.Range("K11").Value = 1.037049963460744
$ws.Range("L11").Value = 1.02837953558615
$ws.Range("M11").Value = 1.043179533237487
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023922207527438
$ws.Range("D12").Value = 1.033279294753137
$ws.Range("E12").Value = 1.024559715216829
$ws.Range("F12").Value = 1.039384930180451
$ws.Range("I12").Value = 1.031316014321759
$ws.Range("J12").Value = 1.030657404064482
$ws.Range("K12").Value = 1.036901624060661
$ws.Range("L12").Value = 1.028215306461055
$ws.Range("M12").Value = 1.042984336630109
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023975488449999
$ws.Range("D13").Value = 1.033320944107484
$ws.Range("E13").Value = 1.024604879171865
$ws.Range("F13").Value = 1.039436596132488
$ws.Range("I13").Value = 1.031325017317561
$ws.Range("J13").Value = 1.03069180359586
$ws.Range("K13").Value = 1.036933447355234
$ws.Range("L13").Value = 1.02825053315011
$ws.Range("M13").Value = 1.043026206513031
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024150081309047
$ws.Range("D14").Value = 1.033457417158671
$ws.Range("E14").Value = 1.024752886668094
$ws.Range("F14").Value = 1.039605904862583
$ws.Range("I14").Value = 1.031354474539408
$ws.Range("J14").Value = 1.030804512215123
$ws.Range("K14").Value = 1.037037703479368
$ws.Range("L14").Value = 1.028365959898601
$ws.Range("M14").Value = 1.043163398022536
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024257664410968
$ws.Range("D15").Value = 1.033541507167938
$ws.Range("E15").Value = 1.024844097663895
$ws.Range("F15").Value = 1.039710237923685
$ws.Range("I15").Value = 1.031372592133688
$ws.Range("J15").Value = 1.030873952642788
$ws.Range("K15").Value = 1.037101927421978
$ws.Range("L15").Value = 1.028437081089731
$ws.Range("M15").Value = 1.043247927548153
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024884007326327
$ws.Range("D16").Value = 1.034031014625463
$ws.Range("E16").Value = 1.025375269749178
$ws.Range("F16").Value = 1.040317750570225
$ws.Range("I16").Value = 1.031477555867851
$ws.Range("J16").Value = 1.031278078530242
$ws.Range("K16").Value = 1.037475559338267
$ws.Range("L16").Value = 1.028851081827209
$ws.Range("M16").Value = 1.043739945250767
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025277038734184
$ws.Range("D17").Value = 1.034338128085611
$ws.Range("E17").Value = 1.025708709884793
$ws.Range("F17").Value = 1.040699045846996
$ws.Range("I17").Value = 1.031542966266989
$ws.Range("J17").Value = 1.031531533727101
$ws.Range("K17").Value = 1.037709770307701
$ws.Range("L17").Value = 1.029110812293812
$ws.Range("M17").Value = 1.044048591280732
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025506338732818
$ws.Range("D18").Value = 1.034517282858681
$ws.Range("E18").Value = 1.025903289624391
$ws.Range("F18").Value = 1.040921527416165
$ws.Range("I18").Value = 1.031580963739994
$ws.Range("J18").Value = 1.031679354400668
$ws.Range("K18").Value = 1.037846324243439
$ws.Range("L18").Value = 1.029262322405849
$ws.Range("M18").Value = 1.044228624893824
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02558453281817
$ws.Range("D19").Value = 1.034578373482767
$ws.Range("E19").Value = 1.025969651480271
$ws.Range("F19").Value = 1.040997401164866
$ws.Range("I19").Value = 1.031593893541256
$ws.Range("J19").Value = 1.031729754788909
$ws.Range("K19").Value = 1.03789287585341
$ws.Range("L19").Value = 1.029313985744307
$ws.Range("M19").Value = 1.044290012628002
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025234864861443
$ws.Range("D20").Value = 1.034305175570186
$ws.Range("E20").Value = 1.02567292562346
$ws.Range("F20").Value = 1.040658128347205
$ws.Range("I20").Value = 1.031535964411263
$ws.Range("J20").Value = 1.031504341976893
$ws.Range("K20").Value = 1.037684647614534
$ws.Range("L20").Value = 1.029082944247792
$ws.Range("M20").Value = 1.044015475896612
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024098665838249
$ws.Range("D21").Value = 1.033417228315902
$ws.Range("E21").Value = 1.024709298253575
$ws.Range("F21").Value = 1.039556044250606
$ws.Range("I21").Value = 1.031345806781622
$ws.Range("J21").Value = 1.030771322974828
$ws.Range("K21").Value = 1.037007005084534
$ws.Range("L21").Value = 1.028331968979763
$ws.Range("M21").Value = 1.043122998268205
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023384758382938
$ws.Range("D22").Value = 1.032859134419642
$ws.Range("E22").Value = 1.024104242838181
$ws.Range("F22").Value = 1.038863833970388
$ws.Range("I22").Value = 1.031224850811751
$ws.Range("J22").Value = 1.030310310146395
$ws.Range("K22").Value = 1.036580433000307
$ws.Range("L22").Value = 1.027859929745268
$ws.Range("M22").Value = 1.042561918603052
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02376316993465
$ws.Range("D23").Value = 1.033154971825678
$ws.Range("E23").Value = 1.024424916545274
$ws.Range("F23").Value = 1.03923071975992
$ws.Range("I23").Value = 1.031289104077685
$ws.Range("J23").Value = 1.030554714315767
$ws.Range("K23").Value = 1.036806615138076
$ws.Range("L23").Value = 1.02811015422714
$ws.Range("M23").Value = 1.042859351911092
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025253921260221
$ws.Range("D24").Value = 1.034320065329635
$ws.Range("E24").Value = 1.025689094713596
$ws.Range("F24").Value = 1.040676616961626
$ws.Range("I24").Value = 1.031539128728367
$ws.Range("J24").Value = 1.031516628806128
$ws.Range("K24").Value = 1.037695999655453
$ws.Range("L24").Value = 1.029095536575849
$ws.Range("M24").Value = 1.044030439295977
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026985871177855
$ws.Range("D25").Value = 1.035672899085462
$ws.Range("E25").Value = 1.027159620203238
$ws.Range("F25").Value = 1.042357575457506
$ws.Range("I25").Value = 1.031823174562657
$ws.Range("J25").Value = 1.032632268350838
$ws.Range("K25").Value = 1.038725827864696
$ws.Range("L25").Value = 1.030239553651163
$ws.Range("M25").Value = 1.045389637757078
